# CitySize.xlsx — "Correct extra 0 in CitySize"
#
# The O column built a hex "class id" by concatenating the B column
# (already a zero-padded hex byte like "0x01") with a literal "00" suffix,
# which produced an extra leading zero (e.g. "0x0100" instead of "0x100").
# Fix: strip the single zero that follows "0x" in B before appending "00".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# O2 holds the first (non-shared) copy of the formula.
$ws.Range("O2").Formula = '=SUBSTITUTE(B2,"0x0","0x")&"00"'

# O3:O18 is a shared-formula block anchored at O3.
$ws.Range("O3:O18").Formula = '=SUBSTITUTE(B3,"0x0","0x")&"00"'

# Column B keeps its existing "best fit" width; column C is widened to 7
# characters (it had been sharing B's bestFit width before).
$ws.Columns.Item(3).ColumnWidth = 6.166666666666667
